# Apply the "quadratic cost function" edit to the 3bus workbook.
#
# Sheet "gen" gains a second cost-coefficient column ("b") next to the
# renamed linear coefficient column ("a", formerly "cost"); the old
# min/max bound columns shift right to make room and the bound values
# change (min 10 -> 0, max stays 400 but now lives in the new column E).

$wb = $excel.ActiveWorkbook

$wsGen = $wb.Worksheets.Item(1)   # "gen"
$wsLin = $wb.Worksheets.Item(2)   # "lin"
$wsRen = $wb.Worksheets.Item(3)   # "ren"
$wsDem = $wb.Worksheets.Item(4)   # "dem"

# --- "gen" sheet: add quadratic cost coefficient column -------------------
# Old layout:  bus | cost | min | max
# New layout:  bus |  a   |  b  | min | max
# The old "min"/"max" columns (C, D) shift right to D/E; column D's old
# "max" values move into the new column E, and column D is repurposed to
# hold the (now explicit) "min" bound.

# Preserve the old "max" values (column D) into the new column E before
# column D is repurposed to hold "min".
$wsGen.Range("E2").Value = 400
$wsGen.Range("E3").Value = 400

# Column D becomes the "min" bound (was blank/implicit, now explicit 0).
$wsGen.Range("D2").Value = 0
$wsGen.Range("D3").Value = 0

# Headers: rename "cost" -> "a", add "b" in column C (pushing "min"/"max"
# headers out to D1/E1).
$wsGen.Range("B1").Value = "a"
$wsGen.Range("D1").Value = "min"
$wsGen.Range("E1").Value = "max"
$wsGen.Range("C1").Value = "b"

# Update the linear coefficients (column B) and the quadratic coefficients
# now living in column C.
$wsGen.Range("B2").Value = 0.1
$wsGen.Range("B3").Value = 0.1
$wsGen.Range("C3").Value = 20

# Match formatting of the surrounding header/data cells (centered, like
# the rest of the sheet) for the newly written cells.
$wsGen.Range("D2:E3").HorizontalAlignment = -4108
$wsGen.Range("C1:E1").HorizontalAlignment = -4108

# --- selection / active-sheet bookkeeping ---------------------------------
# Reproduce the final selection state recorded in each sheet, finishing on
# "gen" (E3) so it ends up the active tab/sheet of the workbook.
$wsDem.Activate() | Out-Null
$wsDem.Range("A1").Select() | Out-Null

$wsLin.Activate() | Out-Null
$wsLin.Range("A1").Select() | Out-Null

$wsRen.Activate() | Out-Null
$wsRen.Range("A1").Select() | Out-Null

$wsGen.Activate() | Out-Null
$wsGen.Range("E3").Select() | Out-Null
